# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-15 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-12-16 Monday", 2)

# The worksheet table alternates a row of 5 problems with 3 blank rows.
# Update each problem cell directly by (row, column) so that the
# "457÷2=" value -- which moves from column 4 to column 1 of row 9 while
# a *different* value ("578÷2=") takes its old slot -- cannot collide
# with a text-based Find/Replace pass.
$tbl = $d.Tables.Item(1)

# Row 1
$tbl.Cell(1, 1).Range.Text = "162÷9="
$tbl.Cell(1, 2).Range.Text = "143÷4="
$tbl.Cell(1, 3).Range.Text = "224÷8="
$tbl.Cell(1, 4).Range.Text = "484÷8="
$tbl.Cell(1, 5).Range.Text = "336÷9="

# Row 5
$tbl.Cell(5, 1).Range.Text = "698÷5="
$tbl.Cell(5, 2).Range.Text = "917÷3="
$tbl.Cell(5, 3).Range.Text = "929÷8="
$tbl.Cell(5, 4).Range.Text = "340÷6="
$tbl.Cell(5, 5).Range.Text = "855÷3="

# Row 9 (cell contents effectively shift left by one position, with two
# brand-new values appended)
$tbl.Cell(9, 1).Range.Text = "457÷2="
$tbl.Cell(9, 2).Range.Text = "445÷9="
$tbl.Cell(9, 3).Range.Text = "509÷4="
$tbl.Cell(9, 4).Range.Text = "578÷2="
$tbl.Cell(9, 5).Range.Text = "977÷6="

# Row 13
$tbl.Cell(13, 1).Range.Text = "613÷5="
$tbl.Cell(13, 2).Range.Text = "684÷8="
$tbl.Cell(13, 3).Range.Text = "848÷3="
$tbl.Cell(13, 4).Range.Text = "345÷7="
$tbl.Cell(13, 5).Range.Text = "500÷8="

# Row 17
$tbl.Cell(17, 1).Range.Text = "976÷8="
$tbl.Cell(17, 2).Range.Text = "891÷8="
$tbl.Cell(17, 3).Range.Text = "472÷5="
$tbl.Cell(17, 4).Range.Text = "245÷4="
$tbl.Cell(17, 5).Range.Text = "339÷6="
